$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D113").Value = 0.791995474
$ws.Range("D114").Value = 0.788120887
$ws.Range("D115").Value = 0.597740902
$ws.Range("D116").Value = 0.620527487
$ws.Range("C117").Value = 0.241887844
$ws.Range("C118").Value = 0.331651578
$ws.Range("C119").Value = 0.154182215
$ws.Range("C120").Value = 0.166899468
$ws.Range("C121").Value = 0.042359665
$ws.Range("C122").Value = 0.266698307
